$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move values from column D to column E for rows 544, 554, 583, 584 ---
$v544 = $ws.Range("D544").Value()
$ws.Range("E544").Value = $v544
$ws.Range("D544").ClearContents() | Out-Null
$v554 = $ws.Range("D554").Value()
$ws.Range("E554").Value = $v554
$ws.Range("D554").ClearContents() | Out-Null
$v583 = $ws.Range("D583").Value()
$ws.Range("E583").Value = $v583
$ws.Range("D583").ClearContents() | Out-Null
$v584 = $ws.Range("D584").Value()
$ws.Range("E584").Value = $v584
$ws.Range("D584").ClearContents() | Out-Null

# --- Prime shared-string table so the 20 brand-new strings are registered in the
#     exact order the source workbook used (a harmless faraway scratch cell is
#     used and cleared afterwards so it leaves no visible trace). ---
$scratch = $ws.Range("Z1000")
$scratch.Value = "REACTIONS SHOT & FINISHING"
$scratch.Value = "FINISHING TWO TOUCHES COMPETITION"
$scratch.Value = "PROTECT SPACE & MOVEMENT TO FINISHING"
$scratch.Value = "IMPROVING LONG BALL (LEFT FOOT)"
$scratch.Value = "BODY SHAPE & LONG PASS"
$scratch.Value = "PENETRATE RUN & FINISHING"
$scratch.Value = "BUILD UP & IMPROVING LEFT FOOT"
$scratch.Value = "OPPONENT PRESS"
$scratch.Value = "POSITIONING GAME. PLAY FORWARD. DIFERENT BEHAVIORS PLAY ON POCKET"
$scratch.Value = "PATTERNS POSITIONAL.IMPROVING TECHNICAL/TACTICAL & CROSS + FINISIHING"
$scratch.Value = "POST MATCH WBA. AVOID LOSSES"
$scratch.Value = "OPPONENT COVENTRY"
$scratch.Value = "PERSONAL INTERVIEW. U20"
$scratch.Value = "COUNTER ATTACK & FINISHING. 1vs1 (GK)"
$scratch.Value = "COMPETITIVE FINISHING DRILLS"
$scratch.Value = "TRAINING NON SQUAD. OFFENSIVE PATTERNS & FINISHING"
$scratch.Value = "DEFENDER ONE-TWO. BUILD UP CONNECTED WITH 11/9"
$scratch.Value = "DEFENDER ONE-TWO. BUILD UP CONNECTED WITH 7/9"
$scratch.Value = "IMPROVING CROSSES 1st POST & 2nd POST"
$scratch.Value = "ATTACK ON THE BOX & FINISHING"
$scratch.ClearContents() | Out-Null

# --- Apply styles (date format on column A, bold on column B) to new rows 723:797 ---
$ws.Range("A722").Copy() | Out-Null
$ws.Range("A723:A797").PasteSpecial(-4122) | Out-Null
$ws.Range("B722").Copy() | Out-Null
$ws.Range("B723:B797").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Populate new training log rows 723-797 ---
$ws.Range("A723").Value = 45945
$ws.Range("B723").Value = "Luca Kjerrumgaard"
$ws.Range("C723").Value = "REACTIONS SHOT & FINISHING"
$ws.Range("A724").Value = 45945
$ws.Range("B724").Value = "Vivaldo Semedo"
$ws.Range("C724").Value = "REACTIONS SHOT & FINISHING"
$ws.Range("A725").Value = 45945
$ws.Range("B725").Value = "Mamadou Doumbia"
$ws.Range("C725").Value = "REACTIONS SHOT & FINISHING"
$ws.Range("A726").Value = 45945
$ws.Range("B726").Value = "Moussa Sissoko"
$ws.Range("C726").Value = "FINISHING TWO TOUCHES COMPETITION"
$ws.Range("A727").Value = 45945
$ws.Range("B727").Value = "Imran Louza"
$ws.Range("C727").Value = "FINISHING TWO TOUCHES COMPETITION"
$ws.Range("A728").Value = 45945
$ws.Range("B728").Value = "Kwadwo Baah"
$ws.Range("C728").Value = "FINISHING TWO TOUCHES COMPETITION"
$ws.Range("A729").Value = 45945
$ws.Range("B729").Value = "Hector Kyprianou"
$ws.Range("C729").Value = "FINISHING TWO TOUCHES COMPETITION"
$ws.Range("A730").Value = 45946
$ws.Range("B730").Value = "Vivaldo Semedo"
$ws.Range("C730").Value = "PROTECT SPACE & MOVEMENT TO FINISHING"
$ws.Range("A731").Value = 45946
$ws.Range("B731").Value = "Mamadou Doumbia"
$ws.Range("C731").Value = "PROTECT SPACE & MOVEMENT TO FINISHING"
$ws.Range("A732").Value = 45946
$ws.Range("B732").Value = "Luca Kjerrumgaard"
$ws.Range("C732").Value = "PROTECT SPACE & MOVEMENT TO FINISHING"
$ws.Range("A733").Value = 45946
$ws.Range("B733").Value = "Nestory Irankunda"
$ws.Range("C733").Value = "PROTECT SPACE & MOVEMENT TO FINISHING"
$ws.Range("A734").Value = 45946
$ws.Range("B734").Value = "Pierre Dwomoh"
$ws.Range("C734").Value = "IMPROVING LONG BALL (LEFT FOOT)"
$ws.Range("A735").Value = 45946
$ws.Range("B735").Value = "Mattie Pollock"
$ws.Range("C735").Value = "IMPROVING LONG BALL (LEFT FOOT)"
$ws.Range("A736").Value = 45947
$ws.Range("B736").Value = "Pierre Dwomoh"
$ws.Range("C736").Value = "BODY SHAPE & LONG PASS"
$ws.Range("A737").Value = 45949
$ws.Range("B737").Value = "Nepalys Mendy"
$ws.Range("C737").Value = "LONG PASS & BODY SHAPE"
$ws.Range("A738").Value = 45949
$ws.Range("B738").Value = "Pierre Dwomoh"
$ws.Range("C738").Value = "LONG PASS & BODY SHAPE"
$ws.Range("A739").Value = 45949
$ws.Range("B739").Value = "Vivaldo Semedo"
$ws.Range("C739").Value = "PENETRATE RUN & FINISHING"
$ws.Range("A740").Value = 45949
$ws.Range("B740").Value = "Mamadou Doumbia"
$ws.Range("C740").Value = "PENETRATE RUN & FINISHING"
$ws.Range("A741").Value = 45949
$ws.Range("B741").Value = "Formose Mendy"
$ws.Range("C741").Value = "IMPROVING CROSSES 1st POST"
$ws.Range("A742").Value = 45949
$ws.Range("B742").Value = "James Morris"
$ws.Range("C742").Value = "IMPROVING CROSSES 1st POST"
$ws.Range("A743").Value = 45950
$ws.Range("B743").Value = "Max Alleyne"
$ws.Range("C743").Value = "BUILD UP & IMPROVING LEFT FOOT"
$ws.Range("A744").Value = 45950
$ws.Range("B744").Value = "Mattie Pollock"
$ws.Range("C744").Value = "BUILD UP & IMPROVING LEFT FOOT"
$ws.Range("A745").Value = 45950
$ws.Range("B745").Value = "James Abankwah"
$ws.Range("C745").Value = "BUILD UP & IMPROVING LEFT FOOT"
$ws.Range("A746").Value = 45950
$ws.Range("B746").Value = "Max Alleyne"
$ws.Range("E746").Value = "OPPONENT PRESS"
$ws.Range("A747").Value = 45951
$ws.Range("B747").Value = "Hector Kyprianou"
$ws.Range("E747").Value = "POSITIONING GAME. PLAY FORWARD. DIFERENT BEHAVIORS PLAY ON POCKET"
$ws.Range("A748").Value = 45952
$ws.Range("B748").Value = "Vivaldo Semedo"
$ws.Range("C748").Value = "PATTERNS POSITIONAL.IMPROVING TECHNICAL/TACTICAL & CROSS + FINISIHING"
$ws.Range("A749").Value = 45952
$ws.Range("B749").Value = "Caleb Wiley"
$ws.Range("C749").Value = "PATTERNS POSITIONAL.IMPROVING TECHNICAL/TACTICAL & CROSS + FINISIHING"
$ws.Range("A750").Value = 45952
$ws.Range("B750").Value = "Nepalys Mendy"
$ws.Range("C750").Value = "PATTERNS POSITIONAL.IMPROVING TECHNICAL/TACTICAL & CROSS + FINISIHING"
$ws.Range("A751").Value = 45952
$ws.Range("B751").Value = "Formose Mendy"
$ws.Range("C751").Value = "PATTERNS POSITIONAL.IMPROVING TECHNICAL/TACTICAL & CROSS + FINISIHING"
$ws.Range("A752").Value = 45952
$ws.Range("B752").Value = "Max Alleyne"
$ws.Range("C752").Value = "PATTERNS POSITIONAL.IMPROVING TECHNICAL/TACTICAL & CROSS + FINISIHING"
$ws.Range("A753").Value = 45953
$ws.Range("B753").Value = "Hector Kyprianou"
$ws.Range("E753").Value = "POST MATCH WBA. AVOID LOSSES"
$ws.Range("A754").Value = 45954
$ws.Range("B754").Value = "Luca Kjerrumgaard"
$ws.Range("E754").Value = "OPPONENT COVENTRY"
$ws.Range("A755").Value = 45959
$ws.Range("B755").Value = "Mamadou Doumbia"
$ws.Range("C755").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A756").Value = 45959
$ws.Range("B756").Value = "Luca Kjerrumgaard"
$ws.Range("C756").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A757").Value = 45959
$ws.Range("B757").Value = "Vivaldo Semedo"
$ws.Range("C757").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A758").Value = 45959
$ws.Range("B758").Value = "Othmane Maamma"
$ws.Range("C758").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A759").Value = 45959
$ws.Range("B759").Value = "Rocco Vata"
$ws.Range("C759").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A760").Value = 45959
$ws.Range("B760").Value = "Tom Ince"
$ws.Range("C760").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A761").Value = 45959
$ws.Range("B761").Value = "Edo Kayembe"
$ws.Range("C761").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A762").Value = 45959
$ws.Range("B762").Value = "Mattie Pollock"
$ws.Range("C762").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A763").Value = 45959
$ws.Range("B763").Value = "Kevin Keben"
$ws.Range("C763").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A764").Value = 45959
$ws.Range("B764").Value = "Kwadwo Baah"
$ws.Range("C764").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A765").Value = 45959
$ws.Range("B765").Value = "Imran Louza"
$ws.Range("C765").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A766").Value = 45959
$ws.Range("B766").Value = "Hector Kyprianou"
$ws.Range("C766").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A767").Value = 45959
$ws.Range("B767").Value = "Pierre Dwomoh"
$ws.Range("C767").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A768").Value = 45959
$ws.Range("B768").Value = "Nepalys Mendy"
$ws.Range("C768").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A769").Value = 45959
$ws.Range("B769").Value = "Formose Mendy"
$ws.Range("C769").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A770").Value = 45959
$ws.Range("B770").Value = "James Abankwah"
$ws.Range("C770").Value = "FINISHING DIFERENT PROFILES. "
$ws.Range("A771").Value = 45959
$ws.Range("B771").Value = "Othmane Maamma"
$ws.Range("D771").Value = "PERSONAL INTERVIEW. U20"
$ws.Range("A772").Value = 45960
$ws.Range("B772").Value = "Mamadou Doumbia"
$ws.Range("C772").Value = "COUNTER ATTACK & FINISHING. 1vs1 (GK)"
$ws.Range("A773").Value = 45960
$ws.Range("B773").Value = "Luca Kjerrumgaard"
$ws.Range("C773").Value = "COUNTER ATTACK & FINISHING. 1vs1 (GK)"
$ws.Range("A774").Value = 45960
$ws.Range("B774").Value = "Othmane Maamma"
$ws.Range("C774").Value = "COUNTER ATTACK & FINISHING. 1vs1 (GK)"
$ws.Range("A775").Value = 45961
$ws.Range("B775").Value = "Max Alleyne"
$ws.Range("C775").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A776").Value = 45961
$ws.Range("B776").Value = "Moussa Sissoko"
$ws.Range("C776").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A777").Value = 45961
$ws.Range("B777").Value = "Tom Ince"
$ws.Range("C777").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A778").Value = 45961
$ws.Range("B778").Value = "Giorgi Chakvetadze"
$ws.Range("C778").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A779").Value = 45961
$ws.Range("B779").Value = "James Morris"
$ws.Range("C779").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A780").Value = 45961
$ws.Range("B780").Value = "James Abankwah"
$ws.Range("C780").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A781").Value = 45961
$ws.Range("B781").Value = "Formose Mendy"
$ws.Range("C781").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A782").Value = 45961
$ws.Range("B782").Value = "Jeremy Petris"
$ws.Range("C782").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A783").Value = 45961
$ws.Range("B783").Value = "Nepalys Mendy"
$ws.Range("C783").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A784").Value = 45961
$ws.Range("B784").Value = "Nestory Irankunda"
$ws.Range("C784").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A785").Value = 45961
$ws.Range("B785").Value = "Edo Kayembe"
$ws.Range("C785").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A786").Value = 45961
$ws.Range("B786").Value = "Vivaldo Semedo"
$ws.Range("C786").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A787").Value = 45961
$ws.Range("B787").Value = "Othmane Maamma"
$ws.Range("C787").Value = "COMPETITIVE FINISHING DRILLS"
$ws.Range("A788").Value = 45962
$ws.Range("B788").Value = "James Morris"
$ws.Range("C788").Value = "TRAINING NON SQUAD. OFFENSIVE PATTERNS & FINISHING"
$ws.Range("A789").Value = 45962
$ws.Range("B789").Value = "Formose Mendy"
$ws.Range("C789").Value = "TRAINING NON SQUAD. OFFENSIVE PATTERNS & FINISHING"
$ws.Range("A790").Value = 45962
$ws.Range("B790").Value = "Nepalys Mendy"
$ws.Range("C790").Value = "TRAINING NON SQUAD. OFFENSIVE PATTERNS & FINISHING"
$ws.Range("A791").Value = 45962
$ws.Range("B791").Value = "James Abankwah"
$ws.Range("C791").Value = "TRAINING NON SQUAD. OFFENSIVE PATTERNS & FINISHING"
$ws.Range("A792").Value = 45964
$ws.Range("B792").Value = "Max Alleyne"
$ws.Range("C792").Value = "DEFENDER ONE-TWO. BUILD UP CONNECTED WITH 11/9"
$ws.Range("A793").Value = 45964
$ws.Range("B793").Value = "Formose Mendy"
$ws.Range("C793").Value = "DEFENDER ONE-TWO. BUILD UP CONNECTED WITH 7/9"
$ws.Range("A794").Value = 45966
$ws.Range("B794").Value = "Luca Kjerrumgaard"
$ws.Range("C794").Value = "ATTACK ON THE BOX & FINISHING"
$ws.Range("A795").Value = 45966
$ws.Range("B795").Value = "Othmane Maamma"
$ws.Range("C795").Value = "IMPROVING CROSSES 1st POST & 2nd POST"
$ws.Range("A796").Value = 45966
$ws.Range("B796").Value = "Jeremy Petris"
$ws.Range("C796").Value = "IMPROVING CROSSES 1st POST"
$ws.Range("A797").Value = 45966
$ws.Range("B797").Value = "Nestory Irankunda"
$ws.Range("C797").Value = "ATTACK ON THE BOX & FINISHING"
